$wb = $excel.ActiveWorkbook

# ============ Sheet "展览" ============
$ws1 = $wb.Worksheets.Item('展览')

# F-column (想去人数) updates for existing rows above the insertion point
$ws1.Cells.Item(2, 6).Value = 614
$ws1.Cells.Item(3, 6).Value = 5766
$ws1.Cells.Item(4, 6).Value = 68
$ws1.Cells.Item(5, 6).Value = 462
$ws1.Cells.Item(8, 6).Value = 388
$ws1.Cells.Item(9, 6).Value = 1369
$ws1.Cells.Item(12, 6).Value = 3114
$ws1.Cells.Item(13, 6).Value = 1945
$ws1.Cells.Item(18, 6).Value = 166
$ws1.Cells.Item(21, 6).Value = 359
$ws1.Cells.Item(23, 6).Value = 41
$ws1.Cells.Item(24, 6).Value = 3638
$ws1.Cells.Item(25, 6).Value = 1154
$ws1.Cells.Item(26, 6).Value = 2883
$ws1.Cells.Item(27, 6).Value = 287
$ws1.Cells.Item(28, 6).Value = 2242
$ws1.Cells.Item(29, 6).Value = 4179
$ws1.Cells.Item(32, 6).Value = 476
$ws1.Cells.Item(33, 6).Value = 1329
$ws1.Cells.Item(34, 6).Value = 87

# Insert new row at sheet-row 35 (pushes old rows 35-47 down to 36-48)
$ws1.Rows.Item(35).Insert()

# Fill newly inserted row 35 with the new event ("夙七烈")
$ws1.Cells.Item(35, 1).Value = 34
$ws1.Cells.Item(35, 2).NumberFormat = "@"
$ws1.Cells.Item(35, 2).Value = '2024-07-20'
$ws1.Cells.Item(35, 3).Value = '杭州·首届次元格子动漫展嘉宾内场——夙七烈'
$ws1.Cells.Item(35, 4).Value = '钱江世纪城奔竞大道353号 杭州国际博览中心'
$ws1.Cells.Item(35, 5).Value = '2024.07.20 09:00-07.20 17:00'
$ws1.Cells.Item(35, 6).Value = 0
$ws1.Cells.Item(35, 7).Value = 238
$ws1.Cells.Item(35, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=86519'
$ws1.Cells.Item(35, 9).Value = '//i1.hdslb.com/bfs/openplatform/202405/tRrp8oLo1716953763180.jpeg'
$ws1.Cells.Item(35, 1).Style = $ws1.Cells.Item(34, 1).Style

# Renumber index column + update F for old rows 35-47 (now at sheet rows 36-48)
$ws1.Cells.Item(36, 1).Value = 35
$ws1.Cells.Item(36, 6).Value = 13
$ws1.Cells.Item(37, 1).Value = 36
$ws1.Cells.Item(37, 6).Value = 34
$ws1.Cells.Item(38, 1).Value = 37
$ws1.Cells.Item(38, 6).Value = 1021
$ws1.Cells.Item(39, 1).Value = 38
$ws1.Cells.Item(39, 6).Value = 1285
$ws1.Cells.Item(40, 1).Value = 39
$ws1.Cells.Item(40, 6).Value = 70
$ws1.Cells.Item(41, 1).Value = 40
$ws1.Cells.Item(41, 6).Value = 1089
$ws1.Cells.Item(42, 1).Value = 41
$ws1.Cells.Item(42, 6).Value = 704
$ws1.Cells.Item(43, 1).Value = 42
$ws1.Cells.Item(43, 6).Value = 586
$ws1.Cells.Item(44, 1).Value = 43
$ws1.Cells.Item(44, 6).Value = 424
$ws1.Cells.Item(45, 1).Value = 44
$ws1.Cells.Item(45, 6).Value = 11
$ws1.Cells.Item(46, 1).Value = 45
$ws1.Cells.Item(46, 6).Value = 87
$ws1.Cells.Item(47, 1).Value = 46
$ws1.Cells.Item(47, 6).Value = 322
$ws1.Cells.Item(48, 1).Value = 47
$ws1.Cells.Item(48, 6).Value = 3599

# ============ Sheet "演出" ============
$ws2 = $wb.Worksheets.Item('演出')
$ws2.Cells.Item(10, 6).Value = 911

# ============ Sheet "全部类型" ============
$ws4 = $wb.Worksheets.Item('全部类型')
$ws4.Cells.Item(2, 6).Value = 614
$ws4.Cells.Item(3, 6).Value = 5766
$ws4.Cells.Item(4, 6).Value = 68
$ws4.Cells.Item(7, 6).Value = 388
$ws4.Cells.Item(8, 6).Value = 1369
$ws4.Cells.Item(9, 6).Value = 3114
$ws4.Cells.Item(11, 6).Value = 1945
$ws4.Cells.Item(15, 6).Value = 911
$ws4.Cells.Item(17, 6).Value = 166
$ws4.Cells.Item(19, 6).Value = 359
$ws4.Cells.Item(20, 6).Value = 3638
$ws4.Cells.Item(23, 6).Value = 1154
$ws4.Cells.Item(25, 6).Value = 2883
$ws4.Cells.Item(26, 6).Value = 2242
$ws4.Cells.Item(27, 6).Value = 4179
$ws4.Cells.Item(30, 6).Value = 1329
$ws4.Cells.Item(31, 6).Value = 34
$ws4.Cells.Item(34, 6).Value = 1285
$ws4.Cells.Item(35, 6).Value = 70
$ws4.Cells.Item(36, 6).Value = 1089
$ws4.Cells.Item(38, 6).Value = 704
$ws4.Cells.Item(40, 6).Value = 424
$ws4.Cells.Item(43, 6).Value = 11
$ws4.Cells.Item(45, 6).Value = 87
$ws4.Cells.Item(47, 6).Value = 322
$ws4.Cells.Item(48, 6).Value = 3599
